$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Checkboxes in column B that were ticked (value 1) in this revision.
$checkedRows = @(28,29,65,66,67,68,71,72,73,75,76,78,83,84,85,86,89,90,91,92,93,98,99,100,103,104,106,107,110)
foreach ($r in $checkedRows) {
    $ws.Range("B$r").Value = 1
}

# A stray single-space label appended in column D next to the (now empty) row 96.
$ws.Range("D96").Value = " "

# Restore the on-screen selection to match the author's last saved view.
$ws.Range("E110").Select()
